$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 96
$ws.Range("I9").Value = 115
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 115
$ws.Range("L9").Value = 20
$ws.Range("M9").Value = 54
$ws.Range("N9").Value = -358
$ws.Range("H41").Value = 2789.7144
$ws.Range("J41").Value = 2181
$ws.Range("L41").Value = 2181
$ws.Range("N41").Value = -3061
$ws.Range("H46").Value = 10050
$ws.Range("J46").Value = 8100
$ws.Range("L46").Value = 24300
$ws.Range("N46").Value = -24538
$ws.Range("H60").Value = 10050
$ws.Range("J60").Value = 8100
$ws.Range("L60").Value = 24300
$ws.Range("N60").Value = -25268
$ws.Range("H92").Value = 825.6
$ws.Range("I92").Value = 782.25
$ws.Range("K92").Value = 782.25
$ws.Range("M92").Value = 465.75
$ws.Range("J100").Value = 4000
$ws.Range("L100").Value = 4000
$ws.Range("N100").Value = -5082
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = $null
$ws.Range("M111").Value = $null
$ws.Range("N111").Value = $null
$ws.Range("H131").Value = 526.1818
$ws.Range("I131").Value = 526.1818
$ws.Range("K131").Value = 1578.5454
$ws.Range("M131").Value = 3461.4546
$ws.Range("H137").Value = 1566.5
$ws.Range("I137").Value = 1566.5
$ws.Range("K137").Value = 4699.5
$ws.Range("M137").Value = -2149.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 525.2727
$ws.Range("I2").Value = 197.66667
$ws.Range("K2").Value = 197.66667
$ws.Range("M2").Value = -84.66667000000001
$ws.Range("H61").Value = 2249.75
$ws.Range("I61").Value = 2249.75
$ws.Range("K61").Value = 2249.75
$ws.Range("M61").Value = -2037.75
$ws.Range("H102").Value = 2877.25
$ws.Range("I102").Value = 2169.6667
$ws.Range("K102").Value = 2169.6667
$ws.Range("M102").Value = -547.6667000000002
$ws.Range("H112").Value = 14999
$ws.Range("J112").Value = 14999
$ws.Range("L112").Value = 14999
$ws.Range("N112").Value = -17953
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = $null
$ws.Range("N114").Value = $null
$ws.Range("H116").Value = 525.2727
$ws.Range("I116").Value = 197.66667
$ws.Range("K116").Value = 197.66667
$ws.Range("M116").Value = 2096.33333
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("H136").Value = 2249.75
$ws.Range("I136").Value = 2249.75
$ws.Range("K136").Value = 6749.25
$ws.Range("M136").Value = -4199.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 525.2727
$ws.Range("I3").Value = 197.66667
$ws.Range("K3").Value = 197.66667
$ws.Range("M3").Value = -83.66667000000001
$ws.Range("H34").Value = 14000
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -1886
$ws.Range("N34").Value = -20228
$ws.Range("H86").Value = 1767.1666
$ws.Range("I86").Value = 1960.6
$ws.Range("K86").Value = 1960.6
$ws.Range("M86").Value = -837.5999999999999
$ws.Range("H89").Value = 1767.1666
$ws.Range("I89").Value = 1960.6
$ws.Range("K89").Value = 9803
$ws.Range("M89").Value = -4187
$ws.Range("H105").Value = 11840.934
$ws.Range("I105").Value = 12415.286
$ws.Range("K105").Value = 12415.286
$ws.Range("M105").Value = -10668.286
$ws.Range("H135").Value = 98900
$ws.Range("J135").Value = 98900
$ws.Range("L135").Value = 98900
$ws.Range("N135").Value = -109040
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3883.3333
$ws.Range("I105").Value = 5375
$ws.Range("K105").Value = 5375
$ws.Range("M105").Value = -3628
$ws.Range("H132").Value = 7502
$ws.Range("H134").Value = 6637.4165
$ws.Range("I134").Value = 1825
$ws.Range("K134").Value = 5475
$ws.Range("M134").Value = -2940
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 922
$ws.Range("I94").Value = 922
$ws.Range("K94").Value = 2766
$ws.Range("M94").Value = -2090
$ws.Range("H115").Value = 3750
$ws.Range("J115").Value = 3750
$ws.Range("L115").Value = 11250
$ws.Range("N115").Value = -13600
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 29332.666
$ws.Range("I10").Value = 8000
$ws.Range("J10").Value = 39999
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 39999
$ws.Range("M10").Value = -7831
$ws.Range("N10").Value = -40337
$ws.Range("H26").Value = 17666.666
$ws.Range("I26").Value = 6000
$ws.Range("K26").Value = 6000
$ws.Range("M26").Value = -5720
$ws.Range("H50").Value = 17666.666
$ws.Range("I50").Value = 6000
$ws.Range("K50").Value = 6000
$ws.Range("M50").Value = -5502
$ws.Range("H70").Value = 3500
$ws.Range("I70").Value = 3500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 3500
$ws.Range("L70").Value = $null
$ws.Range("N70").Value = $null
$ws.Range("M70").Value = -3230
$ws.Range("H73").Value = 3500
$ws.Range("I73").Value = 3500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 3500
$ws.Range("L73").Value = $null
$ws.Range("N73").Value = $null
$ws.Range("M73").Value = -2564
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2375
$ws.Range("I7").Value = 2375
$ws.Range("K7").Value = 2375
$ws.Range("M7").Value = -2263
$ws.Range("H22").Value = 615
$ws.Range("I22").Value = 525.5
$ws.Range("J22").Value = 794
$ws.Range("K22").Value = 525.5
$ws.Range("L22").Value = 794
$ws.Range("M22").Value = -230.5
$ws.Range("N22").Value = -1384
$ws.Range("H27").Value = 615
$ws.Range("I27").Value = 525.5
$ws.Range("J27").Value = 794
$ws.Range("K27").Value = 525.5
$ws.Range("L27").Value = 794
$ws.Range("M27").Value = -418.5
$ws.Range("N27").Value = -1008
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = $null
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = $null
$ws.Range("N61").Value = $null
$ws.Range("H68").Value = 6000
$ws.Range("I68").Value = 6000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5251
$ws.Range("H71").Value = 6000
$ws.Range("I71").Value = 6000
$ws.Range("K71").Value = 30000
$ws.Range("M71").Value = -26256
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = $null
$ws.Range("N113").Value = $null
$ws.Range("H122").Value = 3103.2222
$ws.Range("I122").Value = 2561.2856
$ws.Range("K122").Value = 7683.8568
$ws.Range("M122").Value = -5233.8568
$ws.Range("H126").Value = 2375
$ws.Range("I126").Value = 2375
$ws.Range("K126").Value = 7125
$ws.Range("M126").Value = -4655
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = $null
$ws.Range("H136").Value = 35099.8
$ws.Range("I136").Value = 32624.75
$ws.Range("K136").Value = 97874.25
$ws.Range("M136").Value = -95324.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10166.667
$ws.Range("I4").Value = 18000
$ws.Range("J4").Value = 6250
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 6250
$ws.Range("M4").Value = -17887
$ws.Range("N4").Value = -6476
$ws.Range("H81").Value = 2399.4
$ws.Range("I81").Value = 2399.4
$ws.Range("K81").Value = 4798.8
$ws.Range("M81").Value = -3737.8
$ws.Range("H84").Value = 2399.4
$ws.Range("I84").Value = 2399.4
$ws.Range("K84").Value = 23994
$ws.Range("M84").Value = -18690
$ws.Range("H96").Value = 595
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").Value = $null
$ws.Range("H132").Value = 1667.8
$ws.Range("I132").Value = 1667.8
$ws.Range("K132").Value = 5003.4
$ws.Range("M132").Value = -2473.4
$ws.Range("H136").Value = 2799.1667
$ws.Range("I136").Value = 1513.2858
$ws.Range("J136").Value = 4599.4
$ws.Range("K136").Value = 4539.857400000001
$ws.Range("L136").Value = 13798.2
$ws.Range("M136").Value = -1989.857400000001
$ws.Range("N136").Value = -18898.2
